$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4155.5
$ws.Range("I113").Value = 2463.3333
$ws.Range("J113").Value = 5424.625
$ws.Range("K113").Value = 2463.3333
$ws.Range("L113").Value = 5424.625
$ws.Range("M113").Value = 790.6667000000002
$ws.Range("N113").Value = -11932.625
$ws.Range("H123").Value = 49022.727
$ws.Range("J123").Value = 49022.727
$ws.Range("L123").Value = 49022.727
$ws.Range("N123").Value = -58822.727
$ws.Range("H128").Value = 47500
$ws.Range("J128").Value = 47500
$ws.Range("L128").Value = 47500
$ws.Range("N128").Value = -57460
$ws.Range("H129").Value = 857.46155
$ws.Range("I129").Value = 557.8333
$ws.Range("J129").Value = 1114.2858
$ws.Range("K129").Value = 1673.4999
$ws.Range("L129").Value = 3342.8574
$ws.Range("M129").Value = 3326.5001
$ws.Range("N129").Value = -13342.8574
$ws.Range("H130").Value = 36396
$ws.Range("J130").Value = 36396
$ws.Range("L130").Value = 36396
$ws.Range("N130").Value = -46436
$ws.Range("H133").Value = 54126.168
$ws.Range("J133").Value = 54126.168
$ws.Range("L133").Value = 54126.168
$ws.Range("N133").Value = -64246.168

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23532.652
$ws.Range("I32").Value = 4940.016
$ws.Range("J32").Value = 138807
$ws.Range("K32").Value = 4940.016
$ws.Range("L32").Value = 138807
$ws.Range("M32").Value = -4653.016
$ws.Range("N32").Value = -139381
$ws.Range("H45").Value = 4972.357
$ws.Range("I45").Value = 5302.2
$ws.Range("J45").Value = 4789.1113
$ws.Range("K45").Value = 5302.2
$ws.Range("L45").Value = 4789.1113
$ws.Range("M45").Value = -4925.2
$ws.Range("N45").Value = -5543.1113
$ws.Range("H80").Value = 22171.143
$ws.Range("J80").Value = 22171.143
$ws.Range("L80").Value = 22171.143
$ws.Range("N80").Value = -24167.143
$ws.Range("H83").Value = 22171.143
$ws.Range("J83").Value = 22171.143
$ws.Range("L83").Value = 66513.429
$ws.Range("N83").Value = -76497.429
$ws.Range("H103").Value = 35848
$ws.Range("J103").Value = 35848
$ws.Range("L103").Value = 35848
$ws.Range("N103").Value = -38192
$ws.Range("H109").Value = 23958.666
$ws.Range("J109").Value = 23958.666
$ws.Range("L109").Value = 23958.666
$ws.Range("N109").Value = -26732.666
$ws.Range("H131").Value = 38919.57
$ws.Range("J131").Value = 38919.57
$ws.Range("L131").Value = 38919.57
$ws.Range("N131").Value = -48999.57
$ws.Range("H132").Value = 1639.8223
$ws.Range("I132").Value = 1469.85
$ws.Range("J132").Value = 2999.6
$ws.Range("K132").Value = 4409.549999999999
$ws.Range("L132").Value = 8998.799999999999
$ws.Range("M132").Value = -1879.549999999999
$ws.Range("N132").Value = -14058.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 23865.227
$ws.Range("J82").Value = 30062.688
$ws.Range("L82").Value = 30062.688
$ws.Range("N82").Value = -30828.688
$ws.Range("H85").Value = 23865.227
$ws.Range("J85").Value = 30062.688
$ws.Range("L85").Value = 30062.688
$ws.Range("N85").Value = -32714.688
$ws.Range("H122").Value = 29133.846
$ws.Range("J122").Value = 29133.846
$ws.Range("L122").Value = 29133.846
$ws.Range("N122").Value = -38933.84600000001
$ws.Range("H125").Value = 50142.5
$ws.Range("J125").Value = 50142.5
$ws.Range("L125").Value = 50142.5
$ws.Range("N125").Value = -59982.5
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H130").Value = 35902.5
$ws.Range("J130").Value = 35902.5
$ws.Range("L130").Value = 35902.5
$ws.Range("N130").Value = -45942.5
$ws.Range("H135").Value = 67924
$ws.Range("J135").Value = 67924
$ws.Range("L135").Value = 67924
$ws.Range("N135").Value = -78064

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 48869.75
$ws.Range("J20").Value = 48869.75
$ws.Range("L20").Value = 48869.75
$ws.Range("N20").Value = -49341.75
$ws.Range("H30").Value = 48869.75
$ws.Range("J30").Value = 48869.75
$ws.Range("L30").Value = 48869.75
$ws.Range("N30").Value = -49051.75
$ws.Range("H31").Value = 2162.0889
$ws.Range("I31").Value = 1456.2
$ws.Range("J31").Value = 3573.8667
$ws.Range("K31").Value = 1456.2
$ws.Range("L31").Value = 3573.8667
$ws.Range("M31").Value = -1161.2
$ws.Range("N31").Value = -4163.8667
$ws.Range("H34").Value = 2162.0889
$ws.Range("I34").Value = 1456.2
$ws.Range("J34").Value = 3573.8667
$ws.Range("K34").Value = 1456.2
$ws.Range("L34").Value = 3573.8667
$ws.Range("M34").Value = -1254.2
$ws.Range("N34").Value = -3977.8667
$ws.Range("H50").Value = 8982.286
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H109").Value = 10957.143
$ws.Range("J109").Value = 10957.143
$ws.Range("L109").Value = 10957.143
$ws.Range("N109").Value = -13037.143
$ws.Range("H124").Value = 24652.666
$ws.Range("J124").Value = 24652.666
$ws.Range("L124").Value = 24652.666
$ws.Range("N124").Value = -29562.666
$ws.Range("H128").Value = 48869.75
$ws.Range("J128").Value = 48869.75
$ws.Range("L128").Value = 48869.75
$ws.Range("N128").Value = -58829.75
$ws.Range("H135").Value = 53579.855
$ws.Range("J135").Value = 53579.855
$ws.Range("L135").Value = 53579.855
$ws.Range("N135").Value = -63719.855

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 9810.909
$ws.Range("J93").Value = 9810.909
$ws.Range("L93").Value = 9810.909
$ws.Range("N93").Value = -13554.909
$ws.Range("H123").Value = 34575.332
$ws.Range("J123").Value = 34575.332
$ws.Range("L123").Value = 34575.332
$ws.Range("N123").Value = -39475.332

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2710.7778
$ws.Range("I7").Value = 2549.625
$ws.Range("K7").Value = 2549.625
$ws.Range("M7").Value = -2437.625
$ws.Range("H63").Value = 26266.428
$ws.Range("J63").Value = 26266.428
$ws.Range("L63").Value = 26266.428
$ws.Range("N63").Value = -27764.428
$ws.Range("H66").Value = 26266.428
$ws.Range("J66").Value = 26266.428
$ws.Range("L66").Value = 78799.284
$ws.Range("N66").Value = -86287.284
$ws.Range("H126").Value = 2710.7778
$ws.Range("I126").Value = 2549.625
$ws.Range("K126").Value = 7648.875
$ws.Range("M126").Value = -5178.875
$ws.Range("H127").Value = 55590
$ws.Range("J127").Value = 55590
$ws.Range("L127").Value = 55590
$ws.Range("N127").Value = -65510
$ws.Range("H128").Value = 52329.832
$ws.Range("J128").Value = 52329.832
$ws.Range("L128").Value = 52329.832
$ws.Range("N128").Value = -62289.832
$ws.Range("H129").Value = 45282.25
$ws.Range("J129").Value = 45282.25
$ws.Range("L129").Value = 45282.25
$ws.Range("N129").Value = -55282.25
$ws.Range("H130").Value = 55903.332
$ws.Range("J130").Value = 55903.332
$ws.Range("L130").Value = 55903.332
$ws.Range("N130").Value = -65943.33199999999
$ws.Range("H131").Value = 25930
$ws.Range("J131").Value = 25930
$ws.Range("L131").Value = 25930
$ws.Range("N131").Value = -36010

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29882.857
$ws.Range("J64").Value = 29882.857
$ws.Range("L64").Value = 29882.857
$ws.Range("N64").Value = -30378.857
$ws.Range("H67").Value = 29882.857
$ws.Range("J67").Value = 29882.857
$ws.Range("L67").Value = 29882.857
$ws.Range("N67").Value = -31598.857
$ws.Range("H75").Value = 29993.111
$ws.Range("J75").Value = 29993.111
$ws.Range("L75").Value = 29993.111
$ws.Range("N75").Value = -31865.111
$ws.Range("H78").Value = 29993.111
$ws.Range("J78").Value = 29993.111
$ws.Range("L78").Value = 89979.333
$ws.Range("N78").Value = -99339.333
$ws.Range("H109").Value = 29738.5
$ws.Range("J109").Value = 29738.5
$ws.Range("L109").Value = 29738.5
$ws.Range("N109").Value = -32512.5
$ws.Range("H127").Value = 19219.75
$ws.Range("J127").Value = 19219.75
$ws.Range("L127").Value = 19219.75
$ws.Range("N127").Value = -29139.75
$ws.Range("H128").Value = 47546.363
$ws.Range("J128").Value = 47546.363
$ws.Range("L128").Value = 47546.363
$ws.Range("N128").Value = -57506.363
$ws.Range("H138").Value = 49156.5
$ws.Range("L138").Value = 49156.5
$ws.Range("N138").Value = -49156.5
